# Timesheet update: advance melatonin sleep-efficiency / summary data
# (see commit: "Update and advance melatonin sleep efficiency and summary data")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append a new line of work to the D15 "Work carried out" note.
$d15 = $ws.Range("D15")
$newNote = "Data extracted from new melatonin pdfs" + "`n" + `
    "Split melatonin AE table: > 100 patients | > 3 mths" + "`n" + `
    "Data extracted from new prodrome pdfs" + "`n" + `
    "Remove duplicate references from EndNote library"
$d15.Value = $newNote

# 2) That extra line makes the wrapped-text row taller (57 -> 76pt).
$ws.Rows.Item(15).RowHeight = 76

# 3) Hours worked that week increased from 4 to 11 (formulas below ripple
#    through automatically: TOTAL HOURS, the paid/worked difference, and
#    the credit/debit running balance).
$ws.Range("E15").Value = 11

# 4) Scroll the saved view down so row 13 is at the top (user had scrolled
#    the window while the selection stayed on F15).
$excel.ActiveWindow.ScrollRow = 13
